$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated LR-pairs data: the sending-cluster set now also includes "ECs"
# (following Dr Hou's advice), so the table grows from 2 sending clusters x
# 3 target clusters (6 rows) to 3 sending clusters x 3 target clusters (9 rows).
$data = @(
    @("ECs","Fgf2","Fgfrl1","ECs",2,0.6666666666666666,0.6462393333333333,1.938718,0.03461850536298827,0.03461850536298827,3,1,0.8052596666666667,2.415779,0.121346445909118,0.121346445909118,0.5203904701468889,4.683514231322,0.004200832588484367,0.004200832588484365),
    @("ECs","Fgf2","Fgfrl1","FAPs",2,0.6666666666666666,0.6462393333333333,1.938718,0.03461850536298827,0.03461850536298827,3,1,3.018163,9.054489,0.4548139791235885,0.4548139791235885,1.950455645011333,17.554100805102,0.01574498017545198,0.01574498017545198),
    @("ECs","Fgf2","Fgfrl1","sCs",2,0.6666666666666666,0.6462393333333333,1.938718,0.03461850536298827,0.03461850536298827,3,1,2.812615666666666,8.437847,0.4238395749672935,0.4238395749672934,1.817622873349555,16.358605860146,0.01467269259905192,0.01467269259905191),
    @("FAPs","Fgf2","Fgfrl1","ECs",3,1,15.322826,45.968478,0.8208310864042159,0.8208310864042158,3,1,0.8052596666666667,2.415779,0.121346445909118,0.121346445909118,12.33885375715133,111.049683814362,0.09960493502687172,0.09960493502687169),
    @("FAPs","Fgf2","Fgfrl1","FAPs",3,1,15.322826,45.968478,0.8208310864042159,0.8208310864042158,3,1,3.018163,9.054489,0.4548139791235885,0.4548139791235885,46.246786488638,416.2210783977421,0.3733254525958395,0.3733254525958394),
    @("FAPs","Fgf2","Fgfrl1","sCs",3,1,15.322826,45.968478,0.8208310864042159,0.8208310864042158,3,1,2.812615666666666,8.437847,0.4238395749672935,0.4238395749672934,43.09722046520733,387.874984186866,0.3479006987815046,0.3479006987815045),
    @("sCs","Fgf2","Fgfrl1","ECs",3,1,2.698388,8.095164,0.1445504082327959,0.1445504082327959,3,1,0.8052596666666667,2.415779,0.121346445909118,0.121346445909118,2.172903021417333,19.556127192756,0.01754067829376189,0.01754067829376188),
    @("sCs","Fgf2","Fgfrl1","FAPs",3,1,2.698388,8.095164,0.1445504082327959,0.1445504082327959,3,1,3.018163,9.054489,0.4548139791235885,0.4548139791235885,8.144174821244,73.29757339119601,0.06574354635229702,0.06574354635229702),
    @("sCs","Fgf2","Fgfrl1","sCs",3,1,2.698388,8.095164,0.1445504082327959,0.1445504082327959,3,1,2.812615666666666,8.437847,0.4238395749672935,0.4238395749672934,7.589528363545333,68.305755271908,0.06126618358673697,0.06126618358673696)
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $i + 2
    $vals = $data[$i]
    for ($j = 0; $j -lt $vals.Length; $j++) {
        $ws.Cells.Item($row, $j + 1).Value = $vals[$j]
    }
}
